$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '66.498.26'
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").Value = '3.188.45'
$ws.Range("E3").Value = '  +0.91%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.558'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.77%  '
$ws.Range("D9").Value = '3.185.17'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.519'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.16%  '
$ws.Range("D15").Value = '3.707.55'
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("E16").Value = '  +4.73%  '
$ws.Range("D17").Value = '66.441.49'
$ws.Range("E17").Value = '  +2.34%  '
$ws.Range("D18").Value = '3.187.76'
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '520.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.740'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("E23").Value = '  +4.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("E29").Value = '  +6.87%  '
$ws.Range("E30").Value = '  +13.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("E33").Value = '  +2.62%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '511.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.99'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("E40").Value = '  +10.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.44%  '
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.302'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.51%  '
$ws.Range("D44").Value = '0.0₃0673'
$ws.Range("E44").Value = '  +13.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").Value = '2.902.59'
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.59%  '
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.74%  '
